$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before S (discharge_port3), shifting discharge_rate etc. one column right.
$ws.Columns("S:S").Insert()

# Match the new column's width to its neighbours (Q:R already share the 69.66-char width).
$ws.Range("S1").EntireColumn.ColumnWidth = $ws.Range("R1").EntireColumn.ColumnWidth

# Re-assert the total_port_cost formula as one range write so rows 3:4 stay a shared formula
# (same as the original U3+V3 / U4+V4 pair did before the column shift).
$ws.Range("X3:X4").Formula = "=V3+W3"

# Fill in the real discharge port names (replacing the old placeholder text) and the new
# discharge_port3 values, in the same order the author appears to have typed them.
$ws.Range("R2").Value = "TIANJIN"
$ws.Range("R4").Value = "RIZHAO"
$ws.Range("S4").Value = "TIANJIN"

# Label the newly inserted column.
$ws.Range("S1").Value = "discharge_port3"
